# Suppress the "are you sure you want to delete" alert for sheet deletion.
$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Update "My Sheet" row 2 to hold the "sarvi with filter" data (previously
# stored in row 3 of "My Sheet" before the consolidation).
$ws = $wb.Worksheets.Item("My Sheet")

$ws.Range("A2").Value = "2020-10-06 00:00:00"
$ws.Range("B2").Value = "2020-10-07 00:00:00"
$ws.Range("C2").Value = "sarvi with filter"
$ws.Range("D2").Value = "{'created_date': {'`$gte': datetime.datetime(2020, 10, 6, 0, 0), '`$lt': datetime.datetime(2020, 10, 7, 0, 0)}, 'client_id': 'MNRNJVXE', 'function_name': {'`$ne': 'authorize'}, 'user_id': {'`$ne': None}}"
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = "My Sheet"

# Remove the now-redundant row that used to live at row 3.
$ws.Rows.Item(3).Delete()

# Remove the "Main" and "Sheet 1" worksheets entirely - only "My Sheet"
# remains in the cleaned-up workbook.
$wb.Worksheets.Item("Main").Delete()
$wb.Worksheets.Item("Sheet 1").Delete()

$excel.DisplayAlerts = $true
